# MCH154 collection update: add two new series rows (MCH154-1, MCH154-2)
# to the finding-aid worksheet, matching the "Updated MCH102 to MCH251" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MCH154-1 / VARIOUS BOOKS, PAMPHLETS -----------------------------
$ws.Range("A2").Value = "MCH154-1"
$ws.Range("C2").Value = "VARIOUS BOOKS, PAMPHLETS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21O | GRAP COUNT NUMER: NONE"

# --- Row 3: MCH154-2 / BOOKS -------------------------------------------------
$ws.Range("A3").Value = "MCH154-2"
$ws.Range("C3").Value = "BOOKS"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 21O | GRAP COUNT NUMER: NONE"

# Column D / H are part of the filled-in record range but carry no text --
# still pick up the row's formatting like the rest of the row.
$ws.Range("D2:D3").Value = ""
$ws.Range("H2:H3").Value = ""

# Apply the record styling (10pt Calibri, automatic/theme text colour) to the
# new rows. Column B (alternativeIdentifiers) is intentionally left
# untouched/unstyled, matching the source data which has no value there.
# (Union/comma ranges aren't reliable here, so style each contiguous block
# separately -- the engine dedupes identical resulting styles anyway.)
$ws.Range("A2:A3").Font.Name = "Calibri"
$ws.Range("A2:A3").Font.ThemeColor = 1

$ws.Range("C2:E3").Font.Name = "Calibri"
$ws.Range("C2:E3").Font.ThemeColor = 1

$ws.Range("F2:F3").Font.Name = "Calibri"
$ws.Range("F2:F3").Font.ThemeColor = 1

$ws.Range("G2:H3").Font.Name = "Calibri"
$ws.Range("G2:H3").Font.ThemeColor = 1
